# Refitting NCDEs to individual patients (for manuscript figure)
# Adds a "Label" column (H) marking Control (0) vs MDD (1) rows, and
# updates a handful of refit prediction/error values in the first
# (100-iteration) block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "Label" header in H1, matching the style of the other headers ---
$ws.Range("H1").Value = "Label"
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- Populate the Label column: 0 for Control rows, 1 for MDD rows ---
$labels = @{
    2  = 0; 3  = 0; 4  = 0; 5  = 0; 6  = 0
    7  = 1; 8  = 1; 9  = 1; 10 = 1; 11 = 1
    12 = 0; 13 = 0; 14 = 0; 15 = 0; 16 = 0
    17 = 1; 18 = 1; 19 = 1; 20 = 1; 21 = 1
}
foreach ($row in $labels.Keys) {
    $ws.Cells.Item($row, 8).Value = $labels[$row]
}

# --- Update refit Prediction/Error values (100-iteration block) ---
$ws.Range("D2").Value = 0.5042734937336213
$ws.Range("E2").Value = 0.5042734937336213

$ws.Range("D4").Value = 0.2138758249003624
$ws.Range("E4").Value = 0.2138758249003624

$ws.Range("D6").Value = 0.4219766385357361
$ws.Range("E6").Value = 0.4219766385357361

$ws.Range("D7").Value = 0.3994225207747451
$ws.Range("E7").Value = 0.600577479225255

$ws.Range("D10").Value = 0.4242472482996256
$ws.Range("E10").Value = 0.5757527517003744

$ws.Range("D11").Value = 0.3187520320311939
$ws.Range("E11").Value = 0.6812479679688062
